$d = $word.ActiveDocument

# Locate the contractor ("Đơn vị thi công: Công ty TNHH Nghia Van") table:
# it's the first table whose content starts after the "Nghia Van" text.
$findRng = $d.Content
$found = $findRng.Find.Execute("Nghia Van", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$targetIndex = -1
for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $t = $d.Tables.Item($i)
    if ($t.Range.Start -ge $findRng.End) {
        $targetIndex = $i
        break
    }
}

$tbl = $d.Tables.Item($targetIndex)

# Table was "auto" width; make it a fixed 8651 dxa (= 432.55 pt) width.
$tbl.PreferredWidthType = 3
$tbl.PreferredWidth = 432.55

# Widen the 3rd (last) column from 3645 dxa (182.25 pt) to 3933 dxa (196.65 pt).
$tbl = $d.Tables.Item($targetIndex)
$tbl.Columns.Item(3).Width = 196.65
